$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.212.37"
$ws.Range("E2").Value = "  +1.21%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.277.66"
$ws.Range("E3").Value = "  +0.57%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.70"
$ws.Range("E5").Value = "  +1.69%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "186.37"
$ws.Range("E6").Value = "  +3.86%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.136"
$ws.Range("E9").Value = "  +4.36%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.75"
$ws.Range("E10").Value = "  -0.03%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.418"
$ws.Range("E11").Value = "  +0.80%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.845.17"
$ws.Range("E12").Value = "  +0.60%  "

$ws.Range("E13").Value = "  +0.42%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.70"
$ws.Range("E14").Value = "  +1.71%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "68.207.15"
$ws.Range("E15").Value = "  +1.28%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000173"
$ws.Range("E16").Value = "  +2.45%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.269.62"
$ws.Range("E17").Value = "  +0.33%  "

$ws.Range("E18").Value = "  -0.35%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.66"
$ws.Range("E19").Value = "  +1.84%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "382.51"
$ws.Range("E20").Value = "  +1.23%  "

$ws.Range("E21").Value = "  +1.10%  "

$ws.Range("E22").Value = "  +0.10%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.57"
$ws.Range("E23").Value = "  +0.24%  "

$ws.Range("E24").Value = "  +2.10%  "

$ws.Range("E25").Value = "  +0.41%  "

$ws.Range("E27").Value = "  -1.32%  "

$ws.Range("E28").Value = "  -0.10%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.83"
$ws.Range("E29").Value = "  +3.31%  "

$ws.Range("E30").Value = "  +0.66%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.97"
$ws.Range("E31").Value = "  +1.52%  "

$ws.Range("E32").Value = "  +4.62%  "

$ws.Range("E33").Value = "  +0.94%  "

$ws.Range("E35").Value = "  +2.60%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "163.18"
$ws.Range("E36").Value = "  -0.53%  "

$ws.Range("E37").Value = "  +0.22%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.841"
$ws.Range("E38").Value = "  -2.34%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.82"
$ws.Range("E39").Value = "  +1.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.66"
$ws.Range("E40").Value = "  -1.14%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.67"
$ws.Range("E41").Value = "  +0.92%  "

$ws.Range("E42").Value = "  +5.07%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.42"
$ws.Range("E43").Value = "  +2.22%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0693"
$ws.Range("E44").Value = "  +2.55%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.58"
$ws.Range("E45").Value = "  -0.93%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.646.39"
$ws.Range("E46").Value = "  -4.68%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "344.22"
$ws.Range("E47").Value = "  -0.99%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0286"
$ws.Range("E48").Value = "  +1.71%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "32.08"
$ws.Range("E49").Value = "  +3.73%  "

$ws.Range("E51").Value = "  -0.19%  "

